$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 26 for the "Labex" institution type, shifting
# the existing "Fed", "Ntwk", "Pole", "Pltf" and "Site" rows down by one.
$ws.Rows(26).Insert()

$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "Labex"
$ws.Range("C26").Value = "_"
$ws.Range("D26").Value = "Laboratoire d'excellence"

# Renumber the "Level" column for the rows that shifted down.
$ws.Range("A27").Value = 26
$ws.Range("A28").Value = 27
$ws.Range("A29").Value = 28
$ws.Range("A30").Value = 29
$ws.Range("A31").Value = 30

# Reflect the selection left by the author after the edit.
$ws.Range("D28").Select() | Out-Null
